$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# vdW sheet: manual cell edits (values move left by one column across G:J,
# then a new "comment" column J is populated; K keeps its position but its
# XLOOKUP formula now points at I instead of J).
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("vdW")

# Header row
$ws.Range("G1").Value = "kappa (cm3/mol)"
$ws.Range("H1").Value = "epsAD/k (K)"
$ws.Range("I1").Value = "source"
$ws.Range("J1").Value = "comment"

# Row 3
$ws.Range("G3").Value = 20
$ws.Range("H3").Value = 100.2
$ws.Range("I3").Value = 17
$ws.Range("J3").Value = "not realistic, used only for testing/benchmarking, CaseStudy1"
$ws.Range("K3").Formula = "=_xlfn.XLOOKUP(I3,sources!A`$2:A`$40,sources!C`$2:C`$40)"

# Row 4
$ws.Range("G4").Value = 30
$ws.Range("H4").Value = 100.2
$ws.Range("I4").Value = 17
$ws.Range("J4").Value = "not realistic, used only for testing/benchmarking, CaseStudy1"
$ws.Range("K4").Formula = "=_xlfn.XLOOKUP(I4,sources!A`$2:A`$40,sources!C`$2:C`$40)"

# Row 13 (also fixes a typo: id pointer 502 -> 503, which flips the looked-up name to ETHANOL)
$ws.Range("D13").Value = 503
$ws.Range("G13").Value = 25
$ws.Range("H13").Value = 100.2
$ws.Range("I13").Value = 17
$ws.Range("J13").Value = "not realistic, used only for testing/benchmarking, CaseStudy1"
$ws.Range("K13").Formula = "=_xlfn.XLOOKUP(I13,sources!A`$2:A`$40,sources!C`$2:C`$40)"

# Row 14
$ws.Range("G14").Value = 35
$ws.Range("H14").Value = 100.2
$ws.Range("I14").Value = 17
$ws.Range("J14").Value = "not realistic, used only for testing/benchmarking, CaseStudy1"
$ws.Range("K14").Formula = "=_xlfn.XLOOKUP(I14,sources!A`$2:A`$40,sources!C`$2:C`$40)"

$ws.Range("A4").Select()
$ws.Range("E16").Select()

# ---------------------------------------------------------------------------
# PCSAFT sheet: real column insert before M - a "comment" column, pushing the
# old citation column from M to N.
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("PCSAFT")
$ws.Range("M1").EntireColumn.Insert()
$ws.Range("M1").Value = "comment"

$ws.Range("A19").Select()
$ws.Range("M2").Select()

# ---------------------------------------------------------------------------
# CPA sheet: same real column insert before M.
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("CPA")
$ws.Range("M1").EntireColumn.Insert()
$ws.Range("M1").Value = "comment"

$ws.Range("A22").Select()
$ws.Range("A38:XFD38").Select()

# ---------------------------------------------------------------------------
# NRTL sheet: fix a couple of values + typo'd comments.
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("NRTL")
$ws.Range("H5").Value = 0.2
$ws.Range("K5").Value = "vdw, not realistic, just for testing/benchmarking code, CaseStudy1"
$ws.Range("K6").Value = "vdw, not realistic, just for testing/benchmarking code, CaseStudy2"
$ws.Range("K7").Value = "vdw, not realistic, just for testing/benchmarking code, CaseStudy3"

$ws.Range("H5").Select()

# ---------------------------------------------------------------------------
# Cosmetic-only selection changes on a few other sheets.
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("pure")
$ws.Range("A1").Select()
$ws.Range("F13:F15").Select()

$ws = $wb.Worksheets.Item("volume")
$ws.Range("C4:J4").Select()

$ws = $wb.Worksheets.Item("ed")
$ws.Range("A4:E5").Select()

$ws = $wb.Worksheets.Item("ea")
$ws.Range("A4:E5").Select()
